$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# Add "PAGE input: CH4/N2O Shock, 2060" column blocks (cols AQ:AV)
# mirroring the existing 2010/2020/.../2050 shock blocks, for
# both the CH4 table (rows 1-13) and N2O table (rows 15-27).
# ============================================================

# --- Formatting first: clone styles from the neighbouring 2050
#     block (AJ:AO) onto the new AQ:AV block before writing values,
#     so plain numeric writes below keep that style (s=4 header rows,
#     s=5 data rows) instead of Excel minting a fresh General-format style. ---
$ws.Range("AJ3:AO13").Copy()
$ws.Range("AQ3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AJ17:AO27").Copy()
$ws.Range("AQ17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- CH4 table: group header + column headers ---
$ws.Range("AQ2").Value = "PAGE input: CH4 Shock, 2060"
$ws.Range("AQ3").Value = "Year"
$ws.Range("AR3").Value = "IMAGE"
$ws.Range("AS3").Value = "MESSAGE"
$ws.Range("AT3").Value = "MiniCAM"
$ws.Range("AU3").Value = "MERGE"
$ws.Range("AV3").Value = "Policy"

# --- CH4 table: data rows ---
$ws.Range("AQ4").Value = 2010
$ws.Range("AR4").Value = 0
$ws.Range("AS4").Value = 0
$ws.Range("AT4").Value = 0
$ws.Range("AU4").Value = 0
$ws.Range("AV4").Value = 0

$ws.Range("AQ5").Value = 2020
$ws.Range("AR5").Value = 0
$ws.Range("AS5").Value = 0
$ws.Range("AT5").Value = 0
$ws.Range("AU5").Value = 0
$ws.Range("AV5").Value = 0

$ws.Range("AQ6").Value = 2030
$ws.Range("AR6").Value = 0
$ws.Range("AS6").Value = 0
$ws.Range("AT6").Value = 0
$ws.Range("AU6").Value = 0
$ws.Range("AV6").Value = 0

$ws.Range("AQ7").Value = 2040
$ws.Range("AR7").Value = 0
$ws.Range("AS7").Value = 0
$ws.Range("AT7").Value = 0
$ws.Range("AU7").Value = 0
$ws.Range("AV7").Value = 0

$ws.Range("AQ8").Value = 2050
$ws.Range("AR8").Value = 0
$ws.Range("AS8").Value = 0
$ws.Range("AT8").Value = 0
$ws.Range("AU8").Value = 0
$ws.Range("AV8").Value = 0

$ws.Range("AQ9").Value = 2060
$ws.Range("AR9").Value = (7.1155643087639217 / 100000)
$ws.Range("AS9").Value = (7.2944555823906132 / 100000)
$ws.Range("AT9").Value = (7.508513712286557 / 100000)
$ws.Range("AU9").Value = (7.6260473212275143 / 100000)
$ws.Range("AV9").Value = (8.8229397874300955 / 100000)

$ws.Range("AQ10").Value = 2080
$ws.Range("AR10").Value = (1.3418894915573354 / 100000)
$ws.Range("AS10").Value = (1.4055531157619593 / 100000)
$ws.Range("AT10").Value = (1.4544777158365107 / 100000)
$ws.Range("AU10").Value = (1.4165668575349955 / 100000)
$ws.Range("AV10").Value = (1.7364025824845309 / 100000)

$ws.Range("AQ11").Value = 2100
$ws.Range("AR11").Value = (5.4923995536082519 / 10000000)
$ws.Range("AS11").Value = (5.9687074613901105 / 10000000)
$ws.Range("AT11").Value = (6.1714408356783323 / 10000000)
$ws.Range("AU11").Value = (5.7665128834427647 / 10000000)
$ws.Range("AV11").Value = (7.3930841433478277 / 10000000)

$ws.Range("AQ12").Value = 2200
$ws.Range("AR12").Value = (9.0322458401459472 / 100000000000)
$ws.Range("AS12").Value = (9.9343540060914399 / 100000000000)
$ws.Range("AT12").Value = (1.0258888849534743 / 10000000000)
$ws.Range("AU12").Value = (9.4655263627174685 / 100000000000)
$ws.Range("AV12").Value = (1.2292778128752957 / 10000000000)

$ws.Range("AQ13").Value = 2300
$ws.Range("AR13").Value = (1.2501111257279263 / 10000000000000)
$ws.Range("AS13").Value = (1.3788969965844444 / 10000000000000)
$ws.Range("AT13").Value = (1.4233059175694507 / 10000000000000)
$ws.Range("AU13").Value = (1.3145040611561853 / 10000000000000)
$ws.Range("AV13").Value = (1.7053025658242404 / 10000000000000)

# --- N2O table: group header + column headers ---
$ws.Range("AQ16").Value = "PAGE input: N2O Shock, 2060"
$ws.Range("AQ17").Value = "Year"
$ws.Range("AR17").Value = "IMAGE"
$ws.Range("AS17").Value = "MESSAGE"
$ws.Range("AT17").Value = "MiniCAM"
$ws.Range("AU17").Value = "MERGE"
$ws.Range("AV17").Value = "Policy"

# --- N2O table: data rows ---
$ws.Range("AQ18").Value = 2010
$ws.Range("AR18").Value = 0
$ws.Range("AS18").Value = 0
$ws.Range("AT18").Value = 0
$ws.Range("AU18").Value = 0
$ws.Range("AV18").Value = 0

$ws.Range("AQ19").Value = 2020
$ws.Range("AR19").Value = 0
$ws.Range("AS19").Value = 0
$ws.Range("AT19").Value = 0
$ws.Range("AU19").Value = 0
$ws.Range("AV19").Value = 0

$ws.Range("AQ20").Value = 2030
$ws.Range("AR20").Value = 0
$ws.Range("AS20").Value = 0
$ws.Range("AT20").Value = 0
$ws.Range("AU20").Value = 0
$ws.Range("AV20").Value = 0

$ws.Range("AQ21").Value = 2040
$ws.Range("AR21").Value = 0
$ws.Range("AS21").Value = 0
$ws.Range("AT21").Value = 0
$ws.Range("AU21").Value = 0
$ws.Range("AV21").Value = 0

$ws.Range("AQ22").Value = 2050
$ws.Range("AR22").Value = 0
$ws.Range("AS22").Value = 0
$ws.Range("AT22").Value = 0
$ws.Range("AU22").Value = 0
$ws.Range("AV22").Value = 0

$ws.Range("AQ23").Value = 2060
$ws.Range("AR23").Value = (3.4756994042368374 / 10000)
$ws.Range("AS23").Value = (3.2965924978970161 / 10000)
$ws.Range("AT23").Value = (3.2867020376360357 / 10000)
$ws.Range("AU23").Value = (3.3633103433914903 / 10000)
$ws.Range("AV23").Value = (3.4154397442547058 / 10000)

$ws.Range("AQ24").Value = 2080
$ws.Range("AR24").Value = (3.0442179878799878 / 10000)
$ws.Range("AS24").Value = (2.8776869401197781 / 10000)
$ws.Range("AT24").Value = (2.8413111429783646 / 10000)
$ws.Range("AU24").Value = (2.9299055979585542 / 10000)
$ws.Range("AV24").Value = (2.9949878678081908 / 10000)

$ws.Range("AQ25").Value = 2100
$ws.Range("AR25").Value = (1.8109292692546631 / 10000)
$ws.Range("AS25").Value = (1.7303499667413347 / 10000)
$ws.Range("AT25").Value = (1.6545806484028424 / 10000)
$ws.Range("AU25").Value = (1.7257920543678807 / 10000)
$ws.Range("AV25").Value = (1.7937031897590828 / 10000)

$ws.Range("AQ26").Value = 2200
$ws.Range("AR26").Value = (7.3388024830025605 / 100000)
$ws.Range("AS26").Value = (7.101522016244311 / 100000)
$ws.Range("AT26").Value = (6.5696502927544692 / 100000)
$ws.Range("AU26").Value = (6.9260370142663059 / 100000)
$ws.Range("AV26").Value = (7.3215175532643782 / 100000)

$ws.Range("AQ27").Value = 2300
$ws.Range("AR27").Value = (4.5256296990203726 / 100000)
$ws.Range("AS27").Value = (4.3949774506890638 / 100000)
$ws.Range("AT27").Value = (4.0290512913099263 / 100000)
$ws.Range("AU27").Value = (4.2597198992688767 / 100000)
$ws.Range("AV27").Value = (4.5240913177369002 / 100000)

Write-Output "edit complete"